# Update column F (dSF) values for specific rows per repulled data / mean calculation fix.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -4
    10 = 7
    11 = 5
    16 = -2
    21 = -3
    24 = -7
    30 = -2
    31 = -5
    33 = -9
    36 = -1
    37 = -5
    43 = -5
    44 = 0
    47 = 5
    50 = -6
    53 = -1
    55 = -5
    60 = 5
    64 = 0
    65 = -8
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
